$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 key ("quartz") is authored first (establishes shared-string order),
# then the rest of the rows are filled in top-to-bottom, left-to-right.
$ws.Range("A7").Value = "quartz"

$ws.Range("A5").Value = "olivine"
$ws.Range("B5").Value = "Olivine"
$ws.Range("C5").Value = 0.6

$ws.Range("A6").Value = "olivineDesc"
$ws.Range("B6").Value = "It's green."
$ws.Range("C6").Value = 5

$ws.Range("B7").Value = "Quartz"
$ws.Range("C7").Value = 0.6

$ws.Range("A8").Value = "quartzDesc"
$ws.Range("B8").Value = "Colorful and crystally."
$ws.Range("C8").Value = 5

$ws.Range("A9").Value = "peridotite"
$ws.Range("B9").Value = "Peridotite"
$ws.Range("C9").Value = 1

$ws.Range("A10").Value = "peridotiteDesc"
$ws.Range("B10").Value = "Igneous rock."
$ws.Range("C10").Value = 5

$ws.Range("A11").Value = "gabbro"
$ws.Range("B11").Value = "Gabbro"
$ws.Range("C11").Value = 0.6

$ws.Range("A12").Value = "gabbroDesc"
$ws.Range("B12").Value = "Igneous rock."
$ws.Range("C12").Value = 5

$ws.Range("A13").Value = "diorite"
$ws.Range("B13").Value = "Diorite"
$ws.Range("C13").Value = 0.6

$ws.Range("A14").Value = "dioriteDesc"
$ws.Range("B14").Value = "Igneous rock."
$ws.Range("C14").Value = 5

$ws.Range("A15").Value = "granite"
$ws.Range("B15").Value = "Granite"
$ws.Range("C15").Value = 0.6

$ws.Range("A16").Value = "graniteDesc"
$ws.Range("B16").Value = "Igneous rock."
$ws.Range("C16").Value = 5

$ws.Range("A17").Select()
